# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Refresh case numbers for Polonia, Kazajistan, Eslovaquia (no reordering)
# - Filipinas overtakes Noruega & Chequia (rows 43-45 shift down, Filipinas gets new data)
# - Georgia overtakes Guatemala & Sri Lanka (rows 106-108 shift down, Georgia gets new data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 26 de Abril de 2020 a las 10:22"

function Set-CountryRow {
    param($row, $country, $totalCases, $newCases, $activeCases, $recovered, $critical, $deathsToday, $deaths)
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $totalCases
    $ws.Cells.Item($row, 3).Value = $newCases
    $ws.Cells.Item($row, 4).Value = $activeCases
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# --- Row 33: Polonia (refreshed totals, same rank) ---
Set-CountryRow 33 "Polonia" 11395 122 2265 8604 160 2 526

# --- Rows 43-45: Filipinas climbs above Noruega and Chequia ---
Set-CountryRow 43 "Filipinas" 7579 285 862 6216 1 7 501
Set-CountryRow 44 "Noruega" 7499 6 32 7266 50 0 201
Set-CountryRow 45 "Chequia" 7352 0 2471 4662 79 1 219

# --- Row 62: Kazajistan (refreshed totals, same rank) ---
Set-CountryRow 62 "Kazajistan" 2652 51 646 1981 31 0 25

# --- Row 80: Eslovaquia (refreshed totals, same rank) ---
Set-CountryRow 80 "Eslovaquia" 1379 6 394 967 5 1 18

# --- Rows 106-108: Georgia climbs above Guatemala and Sri Lanka ---
Set-CountryRow 106 "Georgia" 485 29 139 340 6 1 6
Set-CountryRow 107 "Guatemala" 473 43 45 415 5 2 13
Set-CountryRow 108 "Sri Lanka" 462 10 118 337 2 0 7
